$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on the numeric-looking columns (G:K) for the new rows
# so values like "32", "9", "106.66" are stored as text, matching the source data
# (the sheet already records these as text with numberStoredAsText ignored-errors).
$ws.Range("G2:K13").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = ' Nov 6 2020'
$ws.Range("B2").Value = ' Abu Dhabi'
$ws.Range("C2").Value = 'Sunrisers won by 6 wickets (with 2 balls remaining)'
$ws.Range("D2").Value = 'Royal Challengers Bangalore'
$ws.Range("E2").Value = 'Sunrisers Hyderabad'
$ws.Range("F2").Value = 'Aaron Finch '
$ws.Range("G2").Value = '32'
$ws.Range("H2").Value = '30'
$ws.Range("I2").Value = '3'
$ws.Range("J2").Value = '1'
$ws.Range("K2").Value = '106.66'

# Row 3
$ws.Range("A3").Value = ' Oct 21 2020'
$ws.Range("B3").Value = ' Abu Dhabi'
$ws.Range("C3").Value = 'RCB won by 8 wickets (with 39 balls remaining)'
$ws.Range("D3").Value = 'Royal Challengers Bangalore'
$ws.Range("E3").Value = 'Kolkata Knight Riders'
$ws.Range("F3").Value = 'Aaron Finch '
$ws.Range("G3").Value = '16'
$ws.Range("H3").Value = '21'
$ws.Range("I3").Value = '2'
$ws.Range("J3").Value = '0'
$ws.Range("K3").Value = '76.19'

# Row 4
$ws.Range("A4").Value = ' Oct 15 2020'
$ws.Range("B4").Value = ' Sharjah'
$ws.Range("C4").Value = 'Kings XI won by 8 wickets'
$ws.Range("D4").Value = 'Royal Challengers Bangalore'
$ws.Range("E4").Value = 'Kings XI Punjab'
$ws.Range("F4").Value = 'Aaron Finch '
$ws.Range("G4").Value = '20'
$ws.Range("H4").Value = '18'
$ws.Range("I4").Value = '2'
$ws.Range("J4").Value = '1'
$ws.Range("K4").Value = '111.11'

# Row 5
$ws.Range("A5").Value = ' Oct 25 2020'
$ws.Range("B5").Value = ' Dubai (DSC)'
$ws.Range("C5").Value = 'Super Kings won by 8 wickets (with 8 balls remaining)'
$ws.Range("D5").Value = 'Royal Challengers Bangalore'
$ws.Range("E5").Value = 'Chennai Super Kings'
$ws.Range("F5").Value = 'Aaron Finch '
$ws.Range("G5").Value = '15'
$ws.Range("H5").Value = '11'
$ws.Range("I5").Value = '3'
$ws.Range("J5").Value = '0'
$ws.Range("K5").Value = '136.36'

# Row 6
$ws.Range("A6").Value = ' Oct 12 2020'
$ws.Range("B6").Value = ' Sharjah'
$ws.Range("C6").Value = 'RCB won by 82 runs'
$ws.Range("D6").Value = 'Royal Challengers Bangalore'
$ws.Range("E6").Value = 'Kolkata Knight Riders'
$ws.Range("F6").Value = 'Aaron Finch '
$ws.Range("G6").Value = '47'
$ws.Range("H6").Value = '37'
$ws.Range("I6").Value = '4'
$ws.Range("J6").Value = '1'
$ws.Range("K6").Value = '127.02'

# Row 7
$ws.Range("A7").Value = ' Oct 17 2020'
$ws.Range("B7").Value = ' Dubai (DSC)'
$ws.Range("C7").Value = 'RCB won by 7 wickets (with 2 balls remaining)'
$ws.Range("D7").Value = 'Royal Challengers Bangalore'
$ws.Range("E7").Value = 'Rajasthan Royals'
$ws.Range("F7").Value = 'Aaron Finch '
$ws.Range("G7").Value = '14'
$ws.Range("H7").Value = '11'
$ws.Range("I7").Value = '0'
$ws.Range("J7").Value = '2'
$ws.Range("K7").Value = '127.27'

# Row 8
$ws.Range("A8").Value = ' Sep 24 2020'
$ws.Range("B8").Value = ' Dubai (DSC)'
$ws.Range("C8").Value = 'Kings XI won by 97 runs'
$ws.Range("D8").Value = 'Royal Challengers Bangalore'
$ws.Range("E8").Value = 'Kings XI Punjab'
$ws.Range("F8").Value = 'Aaron Finch '
$ws.Range("G8").Value = '20'
$ws.Range("H8").Value = '21'
$ws.Range("I8").Value = '3'
$ws.Range("J8").Value = '0'
$ws.Range("K8").Value = '95.23'

# Row 9
$ws.Range("A9").Value = ' Sep 21 2020'
$ws.Range("B9").Value = ' Dubai (DSC)'
$ws.Range("C9").Value = 'RCB won by 10 runs'
$ws.Range("D9").Value = 'Royal Challengers Bangalore'
$ws.Range("E9").Value = 'Sunrisers Hyderabad'
$ws.Range("F9").Value = 'Aaron Finch '
$ws.Range("G9").Value = '29'
$ws.Range("H9").Value = '27'
$ws.Range("I9").Value = '1'
$ws.Range("J9").Value = '2'
$ws.Range("K9").Value = '107.40'

# Row 10
$ws.Range("A10").Value = ' Sep 28 2020'
$ws.Range("B10").Value = ' Dubai (DSC)'
$ws.Range("C10").Value = 'Match tied (RCB won the one-over eliminator)'
$ws.Range("D10").Value = 'Royal Challengers Bangalore'
$ws.Range("E10").Value = 'Mumbai Indians'
$ws.Range("F10").Value = 'Aaron Finch '
$ws.Range("G10").Value = '52'
$ws.Range("H10").Value = '35'
$ws.Range("I10").Value = '7'
$ws.Range("J10").Value = '1'
$ws.Range("K10").Value = '148.57'

# Row 11
$ws.Range("A11").Value = ' Oct 10 2020'
$ws.Range("B11").Value = ' Dubai (DSC)'
$ws.Range("C11").Value = 'RCB won by 37 runs'
$ws.Range("D11").Value = 'Royal Challengers Bangalore'
$ws.Range("E11").Value = 'Chennai Super Kings'
$ws.Range("F11").Value = 'Aaron Finch '
$ws.Range("G11").Value = '2'
$ws.Range("H11").Value = '9'
$ws.Range("I11").Value = '0'
$ws.Range("J11").Value = '0'
$ws.Range("K11").Value = '22.22'

# Row 12
$ws.Range("A12").Value = ' Oct 3 2020'
$ws.Range("B12").Value = ' Abu Dhabi'
$ws.Range("C12").Value = 'RCB won by 8 wickets (with 5 balls remaining)'
$ws.Range("D12").Value = 'Royal Challengers Bangalore'
$ws.Range("E12").Value = 'Rajasthan Royals'
$ws.Range("F12").Value = 'Aaron Finch '
$ws.Range("G12").Value = '8'
$ws.Range("H12").Value = '7'
$ws.Range("I12").Value = '2'
$ws.Range("J12").Value = '0'
$ws.Range("K12").Value = '114.28'

# Row 13
$ws.Range("A13").Value = ' Oct 5 2020'
$ws.Range("B13").Value = ' Dubai (DSC)'
$ws.Range("C13").Value = 'Capitals won by 59 runs'
$ws.Range("D13").Value = 'Royal Challengers Bangalore'
$ws.Range("E13").Value = 'Delhi Capitals'
$ws.Range("F13").Value = 'Aaron Finch '
$ws.Range("G13").Value = '13'
$ws.Range("H13").Value = '14'
$ws.Range("I13").Value = '1'
$ws.Range("J13").Value = '0'
$ws.Range("K13").Value = '92.85'
